$d = $word.ActiveDocument

# Paragraphs that need strikethrough formatting applied (paragraph mark + all runs).
$strikeMarkers = @(
    ": Parallel tasks for:",
    ": Send an order confirmation email.",
    ": Update order status in the CRM system.",
    ": The process ends after all tasks are completed."
)

# Paragraphs that need a yellow highlight applied (paragraph mark + all runs).
$highlightMarkers = @(
    ": Handle errors during the inventory check and order placement.",
    ": Manual intervention in case of errors."
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text

    foreach ($marker in $strikeMarkers) {
        if ($text.Contains($marker)) {
            $p.Range.Select()
            $word.Selection.Font.StrikeThrough = 1
            break
        }
    }

    foreach ($marker in $highlightMarkers) {
        if ($text.Contains($marker)) {
            $p.Range.Select()
            $word.Selection.Font.HighlightColorIndex = 7
            break
        }
    }
}
